# chore: update Sheets via scheduled runner
# Refresh cached market-price / profit figures (columns H-N) for a handful
# of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2787766.5
$ws.Range("I62").Value = 3974973.5
$ws.Range("J62").Value = 17616.5
$ws.Range("K62").Value = 3974973.5
$ws.Range("L62").Value = 17616.5
$ws.Range("M62").Value = -3974349.5
$ws.Range("N62").Value = -18864.5

$ws.Range("H65").Value = 2787766.5
$ws.Range("I65").Value = 3974973.5
$ws.Range("J65").Value = 17616.5
$ws.Range("K65").Value = 19874867.5
$ws.Range("L65").Value = 88082.5
$ws.Range("M65").Value = -19871747.5
$ws.Range("N65").Value = -94322.5

$ws.Range("H86").Value = 1648.1578
$ws.Range("I86").Value = 1099.7
$ws.Range("J86").Value = 2257.5557
$ws.Range("K86").Value = 1099.7
$ws.Range("L86").Value = 2257.5557
$ws.Range("M86").Value = 23.29999999999995
$ws.Range("N86").Value = -4503.5557

$ws.Range("H89").Value = 1648.1578
$ws.Range("I89").Value = 1099.7
$ws.Range("J89").Value = 2257.5557
$ws.Range("K89").Value = 5498.5
$ws.Range("L89").Value = 11287.7785
$ws.Range("M89").Value = 117.5
$ws.Range("N89").Value = -22519.7785

$ws.Range("H141").Value = 3438.5
$ws.Range("I141").Value = 2408.1924
$ws.Range("J141").Value = 7903.1665
$ws.Range("K141").Value = 7224.5772
$ws.Range("L141").Value = 23709.4995
$ws.Range("M141").Value = -2044.5772
$ws.Range("N141").Value = -34069.49950000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 39445.617
$ws.Range("I2").Value = 39445.617
$ws.Range("K2").Value = 39445.617
$ws.Range("M2").Value = -39332.617

$ws.Range("H32").Value = 20248.328
$ws.Range("I32").Value = 4437.14
$ws.Range("K32").Value = 4437.14
$ws.Range("M32").Value = -4150.14

$ws.Range("H97").Value = 10428
$ws.Range("I97").Value = 11497.777
$ws.Range("K97").Value = 11497.777
$ws.Range("M97").Value = -11001.777

$ws.Range("H110").Value = 1392.4445
$ws.Range("I110").Value = 1060.2727
$ws.Range("J110").Value = 1914.4286
$ws.Range("K110").Value = 1060.2727
$ws.Range("L110").Value = 1914.4286
$ws.Range("M110").Value = 984.7273
$ws.Range("N110").Value = -6004.4286

$ws.Range("H116").Value = 39445.617
$ws.Range("I116").Value = 39445.617
$ws.Range("K116").Value = 39445.617
$ws.Range("M116").Value = -37151.617

$ws.Range("H132").Value = 4417.381
$ws.Range("I132").Value = 4001.7144
$ws.Range("J132").Value = 5248.7144
$ws.Range("K132").Value = 12005.1432
$ws.Range("L132").Value = 15746.1432
$ws.Range("M132").Value = -9475.143199999999
$ws.Range("N132").Value = -20806.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 39445.617
$ws.Range("I3").Value = 39445.617
$ws.Range("K3").Value = 39445.617
$ws.Range("M3").Value = -39331.617

$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 1000
$ws.Range("M86").Value = 123

$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 5000
$ws.Range("M89").Value = 616

$ws.Range("H105").Value = 297453.88
$ws.Range("I105").Value = 3185.5
$ws.Range("K105").Value = 3185.5
$ws.Range("M105").Value = -1438.5

$ws.Range("H134").Value = 4091.6858
$ws.Range("I134").Value = 2324.5217
$ws.Range("K134").Value = 6973.5651
$ws.Range("M134").Value = -4438.5651

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1380
$ws.Range("I16").Value = 1603.6666
$ws.Range("K16").Value = 1603.6666
$ws.Range("M16").Value = -1316.6666

$ws.Range("H58").Value = 2284.6875
$ws.Range("I58").Value = 1264.3478
$ws.Range("J58").Value = 4892.222
$ws.Range("K58").Value = 1264.3478
$ws.Range("L58").Value = 4892.222
$ws.Range("M58").Value = -1061.3478
$ws.Range("N58").Value = -5298.222

$ws.Range("H113").Value = 1380
$ws.Range("I113").Value = 1603.6666
$ws.Range("K113").Value = 1603.6666
$ws.Range("M113").Value = 566.3334

$ws.Range("H136").Value = 2284.6875
$ws.Range("I136").Value = 1264.3478
$ws.Range("J136").Value = 4892.222
$ws.Range("K136").Value = 3793.0434
$ws.Range("L136").Value = 14676.666
$ws.Range("M136").Value = -1243.0434
$ws.Range("N136").Value = -19776.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2790.9092
$ws.Range("J34").Value = 2619.8
$ws.Range("L34").Value = 7859.400000000001
$ws.Range("N34").Value = -8027.400000000001

$ws.Range("H39").Value = 8273.549000000001
$ws.Range("I39").Value = 480
$ws.Range("J39").Value = 8533.333000000001
$ws.Range("K39").Value = 1440
$ws.Range("L39").Value = 25599.999
$ws.Range("M39").Value = -1146
$ws.Range("N39").Value = -26187.999

$ws.Range("H55").Value = 3300
$ws.Range("I55").Value = 900
$ws.Range("J55").Value = 3780
$ws.Range("K55").Value = 2700
$ws.Range("L55").Value = 11340
$ws.Range("M55").Value = -2523
$ws.Range("N55").Value = -11694

$ws.Range("H118").Value = 3406
$ws.Range("I118").Value = 1030
$ws.Range("K118").Value = 3090
$ws.Range("M118").Value = -1847

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3021.7
$ws.Range("I80").Value = 2923.077
$ws.Range("J80").Value = 3204.8572
$ws.Range("K80").Value = 2923.077
$ws.Range("L80").Value = 3204.8572
$ws.Range("M80").Value = -1925.077
$ws.Range("N80").Value = -5200.8572

$ws.Range("H83").Value = 3021.7
$ws.Range("I83").Value = 2923.077
$ws.Range("J83").Value = 3204.8572
$ws.Range("K83").Value = 14615.385
$ws.Range("L83").Value = 16024.286
$ws.Range("M83").Value = -9623.385000000002
$ws.Range("N83").Value = -26008.286

$ws.Range("H113").Value = 1340
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 1466.6666
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 1466.6666
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -5806.6666

$ws.Range("H132").Value = 2770.4055
$ws.Range("I132").Value = 2203.6562
$ws.Range("K132").Value = 6610.9686
$ws.Range("M132").Value = -4080.9686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1926.3478
$ws.Range("I132").Value = 1761.9395
$ws.Range("J132").Value = 2343.6924
$ws.Range("K132").Value = 5285.818499999999
$ws.Range("L132").Value = 7031.0772
$ws.Range("M132").Value = -2755.818499999999
$ws.Range("N132").Value = -12091.0772
